# Apply crypto price/volume updates from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "64.659.04"
$cell.Style = "Normal"

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  +4.75%  "
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.095.19"
$cell.Style = "Normal"

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  +3.00%  "
$cell.Style = "Normal"

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.12%  "
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "558.49"
$cell.Style = "Normal"

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  +2.76%  "
$cell.Style = "Normal"

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "144.64"
$cell.Style = "Normal"

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  +9.18%  "
$cell.Style = "Normal"

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -0.19%  "
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "3.089.95"
$cell.Style = "Normal"

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  +3.03%  "
$cell.Style = "Normal"

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.500"
$cell.Style = "Normal"

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  +2.00%  "
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "7.13"
$cell.Style = "Normal"

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  +15.74%  "
$cell.Style = "Normal"

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  +4.48%  "
$cell.Style = "Normal"

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.465"
$cell.Style = "Normal"

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  +4.48%  "
$cell.Style = "Normal"

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000228"
$cell.Style = "Normal"

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  +4.04%  "
$cell.Style = "Normal"

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "35.42"
$cell.Style = "Normal"

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  +2.67%  "
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.595.61"
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "64.620.55"
$cell.Style = "Normal"

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  +4.63%  "
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "3.091.70"
$cell.Style = "Normal"

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  +2.90%  "
$cell.Style = "Normal"

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -0.52%  "
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "6.82"
$cell.Style = "Normal"

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  +3.15%  "
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "483.45"
$cell.Style = "Normal"

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  +0.32%  "
$cell.Style = "Normal"

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.82"
$cell.Style = "Normal"

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  +4.52%  "
$cell.Style = "Normal"

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.678"
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  +1.32%  "
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.57"
$cell.Style = "Normal"

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  +9.02%  "
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "13.40"
$cell.Style = "Normal"

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  +12.31%  "
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "81.01"
$cell.Style = "Normal"

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -1.08%  "
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.78"
$cell.Style = "Normal"

$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  +3.18%  "
$cell.Style = "Normal"

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "8.23"
$cell.Style = "Normal"

$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  +6.63%  "
$cell.Style = "Normal"

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.07"
$cell.Style = "Normal"

$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  +8.63%  "
$cell.Style = "Normal"

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"

$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  +0.02%  "
$cell.Style = "Normal"

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "26.08"
$cell.Style = "Normal"

$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  +1.76%  "
$cell.Style = "Normal"

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  +2.41%  "
$cell.Style = "Normal"

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "2.46"
$cell.Style = "Normal"

$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  +5.74%  "
$cell.Style = "Normal"

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "5.70"
$cell.Style = "Normal"

$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  +0.83%  "
$cell.Style = "Normal"

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "6.24"
$cell.Style = "Normal"

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +6.76%  "
$cell.Style = "Normal"

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "54.85"
$cell.Style = "Normal"

$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = "  -0.36%  "
$cell.Style = "Normal"

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "466.49"
$cell.Style = "Normal"

$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  +5.32%  "
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.0408"
$cell.Style = "Normal"

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  +6.98%  "
$cell.Style = "Normal"

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.0827"
$cell.Style = "Normal"

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  +4.09%  "
$cell.Style = "Normal"

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.88"
$cell.Style = "Normal"

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  +18.92%  "
$cell.Style = "Normal"

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "3.007.33"
$cell.Style = "Normal"

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -4.26%  "
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "8.29"
$cell.Style = "Normal"

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  +2.77%  "
$cell.Style = "Normal"

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -1.21%  "
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "28.40"
$cell.Style = "Normal"

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = "  +7.52%  "
$cell.Style = "Normal"

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.258"
$cell.Style = "Normal"

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  +7.09%  "
$cell.Style = "Normal"

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.11"
$cell.Style = "Normal"

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  +8.36%  "
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  +4.21%  "
$cell.Style = "Normal"

$ws.Range("B49").Value = "PEPE"

$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0$([char]0x2083)0518"
$cell.Style = "Normal"

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  +6.48%  "
$cell.Style = "Normal"

$ws.Range("B50").Value = "Monero"

$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "117.88"
$cell.Style = "Normal"

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  +1.67%  "
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.08"
$cell.Style = "Normal"

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  +2.63%  "
$cell.Style = "Normal"
